# Update "Förändrad" date column (C) for rows 2-18 from serial 45208 (2023-10-09)
# to serial 45212 (2023-10-13), matching the canonical OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 18; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45208) {
        $cell.Value2 = 45212
    }
}
